# Add a new "Data Error Orders" test case row/column to the Input_Data sheet.
#
# Summary of the change (per the target diff):
#   - New shared strings: "data_errors_orders" (A5) and
#     "DataErrorResubmitOrderConfirmationId" (V1, new header column).
#   - New cell A5 = "data_errors_orders" (same style as A2:A4, i.e. copy A4's format).
#   - New cell V1 = "DataErrorResubmitOrderConfirmationId" (same style as the other
#     header cells, i.e. copy U1's format).
#   - Used range grows from A1:U8 to A1:V8.
#   - Column V gets a custom width, and column U's width also grows slightly.
#   - Selection moves to V5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data -------------------------------------------------------------
# Order matters for shared-string insertion order (A5's string must be
# registered before V1's, to match the target sharedStrings.xml ordering).
$ws.Range("A5").Value = "data_errors_orders"
$ws.Range("V1").Value = "DataErrorResubmitOrderConfirmationId"

# --- Formatting: copy styles from the neighboring/previous cells ----------
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Column widths ----------------------------------------------------------
# Target stored widths (OOXML <col width>) are 35.42578125 (U) and
# 35.7109375 (V). The ColumnWidth COM property is quantized to whole pixels
# before being stored, so we pick the closest achievable ColumnWidth values.
$ws.Columns("U").ColumnWidth = 34.583333333333336
$ws.Columns("V").ColumnWidth = 34.75

# --- Selection --------------------------------------------------------------
$ws.Range("V5").Select()
